$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column at N ---------
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14), pushing Late/heading/Outstanding
# (old N/O/P) one column to the right (-> O/P/Q).
$ws.Columns.Item(14).Insert()

# The inserted column should pick up the width of its left neighbour (M),
# matching Excel's normal "insert column" behaviour.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# This sheet becomes the active tab/sheet, with the cursor left on K18.
$ws.Activate()
$null = $ws.Range("K18").Select()
